# Auto-generated edit script applying scheduled market-data refresh
# to the Leve profit sheets (columns H-N) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1695.7142
$ws.Range("J112").Value = 1765.8334
$ws.Range("L112").Value = 5297.5002
$ws.Range("N112").Value = -7513.5002
$ws.Range("H132").Value = 1242.443
$ws.Range("I132").Value = 1255.7164
$ws.Range("J132").Value = 1168.3334
$ws.Range("K132").Value = 3767.1492
$ws.Range("L132").Value = 3505.0002
$ws.Range("M132").Value = -1237.1492
$ws.Range("N132").Value = -8565.0002
$ws.Range("H137").Value = 946.2174
$ws.Range("I137").Value = 1088
$ws.Range("J137").Value = 906.8333
$ws.Range("K137").Value = 3264
$ws.Range("L137").Value = 2720.4999
$ws.Range("M137").Value = -714
$ws.Range("N137").Value = -7820.4999
$ws.Range("H138").Value = 1385.02
$ws.Range("I138").Value = 725.0161000000001
$ws.Range("J138").Value = 2461.8684
$ws.Range("K138").Value = 2175.0483
$ws.Range("L138").Value = 7385.6052
$ws.Range("M138").Value = 2964.9517
$ws.Range("N138").Value = -17665.6052
$ws.Range("H141").Value = 2413.2563
$ws.Range("I141").Value = 768
$ws.Range("J141").Value = 7184.5
$ws.Range("K141").Value = 2304
$ws.Range("L141").Value = 21553.5
$ws.Range("M141").Value = 2876
$ws.Range("N141").Value = -31913.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1510.57
$ws.Range("I32").Value = 1319.6709
$ws.Range("J32").Value = 2228.7144
$ws.Range("K32").Value = 1319.6709
$ws.Range("L32").Value = 2228.7144
$ws.Range("M32").Value = -1032.6709
$ws.Range("N32").Value = -2802.7144
$ws.Range("H61").Value = 872.67566
$ws.Range("I61").Value = 720.65515
$ws.Range("J61").Value = 1423.75
$ws.Range("K61").Value = 720.65515
$ws.Range("L61").Value = 1423.75
$ws.Range("M61").Value = -508.65515
$ws.Range("N61").Value = -1847.75
$ws.Range("H110").Value = 1107.6666
$ws.Range("I110").Value = 1082
$ws.Range("K110").Value = 1082
$ws.Range("M110").Value = 963
$ws.Range("H132").Value = 1645
$ws.Range("I132").Value = 1482.9375
$ws.Range("J132").Value = 2509.3333
$ws.Range("K132").Value = 4448.8125
$ws.Range("L132").Value = 7527.999899999999
$ws.Range("M132").Value = -1918.8125
$ws.Range("N132").Value = -12587.9999
$ws.Range("H135").Value = 41426.5
$ws.Range("J135").Value = 41426.5
$ws.Range("L135").Value = 41426.5
$ws.Range("N135").Value = -51566.5
$ws.Range("H136").Value = 872.67566
$ws.Range("I136").Value = 720.65515
$ws.Range("J136").Value = 1423.75
$ws.Range("K136").Value = 2161.96545
$ws.Range("L136").Value = 4271.25
$ws.Range("M136").Value = 388.0345499999999
$ws.Range("N136").Value = -9371.25
$ws.Range("H139").Value = 79805
$ws.Range("J139").Value = 79805
$ws.Range("L139").Value = 79805
$ws.Range("N139").Value = -90085

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 17223.36
$ws.Range("I134").Value = 1331.5636
$ws.Range("K134").Value = 3994.6908
$ws.Range("M134").Value = -1459.6908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2854
$ws.Range("I99").Value = 2411.25
$ws.Range("J99").Value = 4625
$ws.Range("K99").Value = 2411.25
$ws.Range("L99").Value = 4625
$ws.Range("M99").Value = -913.25
$ws.Range("N99").Value = -7621
$ws.Range("H122").Value = 838
$ws.Range("I122").Value = 800
$ws.Range("J122").Value = 914
$ws.Range("K122").Value = 2400
$ws.Range("L122").Value = 2742
$ws.Range("M122").Value = 50
$ws.Range("N122").Value = -7642
$ws.Range("H126").Value = 2854
$ws.Range("I126").Value = 2411.25
$ws.Range("J126").Value = 4625
$ws.Range("K126").Value = 7233.75
$ws.Range("L126").Value = 13875
$ws.Range("M126").Value = -4763.75
$ws.Range("N126").Value = -18815
$ws.Range("H132").Value = 1319.7161
$ws.Range("I132").Value = 778.6226
$ws.Range("J132").Value = 2343.9285
$ws.Range("K132").Value = 2335.8678
$ws.Range("L132").Value = 7031.7855
$ws.Range("M132").Value = 194.1322
$ws.Range("N132").Value = -12091.7855
$ws.Range("H134").Value = 993.7619
$ws.Range("I134").Value = 899.4865
$ws.Range("K134").Value = 2698.4595
$ws.Range("M134").Value = -163.4594999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 611.7954999999999
$ws.Range("I5").Value = 585.54285
$ws.Range("J5").Value = 713.8889
$ws.Range("K5").Value = 1756.62855
$ws.Range("L5").Value = 2141.6667
$ws.Range("M5").Value = -1644.62855
$ws.Range("N5").Value = -2365.6667
$ws.Range("H122").Value = 550.10254
$ws.Range("I122").Value = 247
$ws.Range("J122").Value = 701.6539
$ws.Range("K122").Value = 2223
$ws.Range("L122").Value = 6314.8851
$ws.Range("M122").Value = 227
$ws.Range("N122").Value = -11214.8851
$ws.Range("H135").Value = 611.7954999999999
$ws.Range("I135").Value = 585.54285
$ws.Range("J135").Value = 713.8889
$ws.Range("K135").Value = 5269.88565
$ws.Range("L135").Value = 6425.0001
$ws.Range("M135").Value = -2734.88565
$ws.Range("N135").Value = -11495.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4570
$ws.Range("I70").Value = 4066.6667
$ws.Range("J70").Value = 4785.7144
$ws.Range("K70").Value = 4066.6667
$ws.Range("L70").Value = 4785.7144
$ws.Range("M70").Value = -3796.6667
$ws.Range("N70").Value = -5325.7144
$ws.Range("H73").Value = 4570
$ws.Range("I73").Value = 4066.6667
$ws.Range("J73").Value = 4785.7144
$ws.Range("K73").Value = 4066.6667
$ws.Range("L73").Value = 4785.7144
$ws.Range("M73").Value = -3130.6667
$ws.Range("N73").Value = -6657.7144
$ws.Range("H102").Value = 1670.75
$ws.Range("I102").Value = 1683.7142
$ws.Range("J102").Value = 1580
$ws.Range("K102").Value = 1683.7142
$ws.Range("L102").Value = 1580
$ws.Range("M102").Value = -61.71419999999989
$ws.Range("N102").Value = -4824
$ws.Range("H132").Value = 2207.9736
$ws.Range("I132").Value = 2088.7407
$ws.Range("J132").Value = 2500.6365
$ws.Range("K132").Value = 6266.222099999999
$ws.Range("L132").Value = 7501.9095
$ws.Range("M132").Value = -3736.222099999999
$ws.Range("N132").Value = -12561.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 506451.34
$ws.Range("J40").Value = 2101
$ws.Range("L40").Value = 2101
$ws.Range("N40").Value = -2373
$ws.Range("H61").Value = 6244.1816
$ws.Range("I61").Value = 7539.5293
$ws.Range("J61").Value = 1840
$ws.Range("K61").Value = 7539.5293
$ws.Range("L61").Value = 1840
$ws.Range("M61").Value = -7337.5293
$ws.Range("N61").Value = -2244
$ws.Range("H113").Value = 6244.1816
$ws.Range("I113").Value = 7539.5293
$ws.Range("J113").Value = 1840
$ws.Range("K113").Value = 7539.5293
$ws.Range("L113").Value = 1840
$ws.Range("M113").Value = -5369.5293
$ws.Range("N113").Value = -6180
$ws.Range("H132").Value = 1689.9412
$ws.Range("I132").Value = 1610.234
$ws.Range("K132").Value = 4830.701999999999
$ws.Range("M132").Value = -2300.701999999999
$ws.Range("H136").Value = 2142.2
$ws.Range("I136").Value = 1283.9678
$ws.Range("J136").Value = 4042.5715
$ws.Range("K136").Value = 3851.9034
$ws.Range("L136").Value = 12127.7145
$ws.Range("M136").Value = -1301.9034
$ws.Range("N136").Value = -17227.7145
$ws.Range("H140").Value = 157344
$ws.Range("J140").Value = 157344
$ws.Range("L140").Value = 157344
$ws.Range("N140").Value = -167704

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 24333.334
$ws.Range("J114").Value = 24333.334
$ws.Range("L114").Value = 24333.334
$ws.Range("N114").Value = -33011.334
$ws.Range("H123").Value = 49820
$ws.Range("J123").Value = 49820
$ws.Range("L123").Value = 49820
$ws.Range("N123").Value = -59620
$ws.Range("H132").Value = 530.8378
$ws.Range("I132").Value = 474.9524
$ws.Range("J132").Value = 850.9091
$ws.Range("K132").Value = 1424.8572
$ws.Range("L132").Value = 2552.7273
$ws.Range("M132").Value = 1105.1428
$ws.Range("N132").Value = -7612.7273
$ws.Range("H136").Value = 797.34784
$ws.Range("I136").Value = 1058.44
$ws.Range("J136").Value = 486.5238
$ws.Range("K136").Value = 3175.32
$ws.Range("L136").Value = 1459.5714
$ws.Range("M136").Value = -625.3200000000002
$ws.Range("N136").Value = -6559.5714
